# #23 Updated/Added new test cases for to-do page features
#
# Four new test-case rows' "Actual Outcome" (col F) and row 4's "Fail/Pass"
# (col G) are updated to reflect new/failed test cases. The new descriptive
# strings are written in the same order the target workbook appended them to
# the shared-string table (F5, F4, F7, F8) so the resulting sharedStrings.xml
# indices line up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("F5").Value = "Item cannot be automatically archived."
$ws.Range("F4").Value = "Item cannot be deleted."
$ws.Range("F7").Value = "Timestamp is not being displayed."
$ws.Range("F8").Value = "The to-do items are not filtered by users. All users see the same set of to-do items."
$ws.Range("G4").Value = "Fail"

# Reflect the author's final scroll position / selection when they saved.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("F8").Select()
